$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all existing data rows
#    (rows 2-497) from 45190 (2023-09-21) to 45192 (2023-09-23).
$ws.Range("C2:C497").Value = 45192

# 2. The previously-last row (497) now gets an explicit row height, matching
#    the pattern produced by the exporter once a row stops being the very
#    last row in the sheet.
$ws.Rows.Item(497).RowHeight = 15

# 3. Append new row 498 - A 44469-2023
$ws.Range("A498").Value = "A 44469-2023"
$ws.Range("B498").Value = 45189
$ws.Range("B498").NumberFormat = "YYYY-MM-DD"
$ws.Range("C498").Value = 45192
$ws.Range("C498").NumberFormat = "YYYY-MM-DD"
$ws.Range("D498").Value = "JÖNKÖPINGS LÄN"
$ws.Range("E498").Value = "GISLAVED"
$ws.Range("G498").Value = 2.3
$ws.Range("H498").Value = 0
$ws.Range("I498").Value = 0
$ws.Range("J498").Value = 0
$ws.Range("K498").Value = 0
$ws.Range("L498").Value = 0
$ws.Range("M498").Value = 0
$ws.Range("N498").Value = 0
$ws.Range("O498").Value = 0
$ws.Range("P498").Value = 0
$ws.Range("Q498").Value = 0
$ws.Range("R498").WrapText = $true

# Row 498 is not the final row anymore (row 499 follows), so it also gets the
# explicit row height treatment.
$ws.Rows.Item(498).RowHeight = 15

# 4. Append new row 499 - A 44742-2023 (the new final row, left without an
#    explicit row height, matching the exporter's behaviour for the last row).
$ws.Range("A499").Value = "A 44742-2023"
$ws.Range("B499").Value = 45190
$ws.Range("B499").NumberFormat = "YYYY-MM-DD"
$ws.Range("C499").Value = 45192
$ws.Range("C499").NumberFormat = "YYYY-MM-DD"
$ws.Range("D499").Value = "JÖNKÖPINGS LÄN"
$ws.Range("E499").Value = "GISLAVED"
$ws.Range("G499").Value = 0.6
$ws.Range("H499").Value = 0
$ws.Range("I499").Value = 0
$ws.Range("J499").Value = 0
$ws.Range("K499").Value = 0
$ws.Range("L499").Value = 0
$ws.Range("M499").Value = 0
$ws.Range("N499").Value = 0
$ws.Range("O499").Value = 0
$ws.Range("P499").Value = 0
$ws.Range("Q499").Value = 0
$ws.Range("R499").WrapText = $true
